function Set-HyperlinkDisplay {
    param($ws, $addr, $text)
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": the fa02844e record is now listed first (row 2), the
# 5ea61590 record moves to row 3 and is reported as newly "Ready for
# handoff" with an updated generation timestamp.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "fa02844e-356c-442b-89d9-021fb9e6330d.md"
$wsOverview.Range("B2").Value = "e2e\fa02844e-356c-442b-89d9-021fb9e6330d.md"
Set-HyperlinkDisplay $wsOverview '$B$2' "e2e\fa02844e-356c-442b-89d9-021fb9e6330d.md"

$wsOverview.Range("A3").Value = "5ea61590-ac6e-4761-8d83-7c8a94ac5562.md"
$wsOverview.Range("B3").Value = "e2e\5ea61590-ac6e-4761-8d83-7c8a94ac5562.md"
Set-HyperlinkDisplay $wsOverview '$B$3' "e2e\5ea61590-ac6e-4761-8d83-7c8a94ac5562.md"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-04 06:52:28"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": same row swap, plus status/date refresh and a stale-handback
# error note for the (now second) 5ea61590 record. Column P widened to fit
# the new error text.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

$wsZhCn.Range("A2").Value = "fa02844e-356c-442b-89d9-021fb9e6330d.md"
Set-HyperlinkDisplay $wsZhCn '$A$2' "fa02844e-356c-442b-89d9-021fb9e6330d.md"
$wsZhCn.Range("G2").Value = "fa02844e-356c-442b-89d9-021fb9e6330d.4bcf33a3301db87c2e28c06185ec766dc2fa2c70.zh-cn.xlf"
$wsZhCn.Range("I2").Value = "fa02844e-356c-442b-89d9-021fb9e6330d.md"
Set-HyperlinkDisplay $wsZhCn '$I$2' "fa02844e-356c-442b-89d9-021fb9e6330d.md"
$wsZhCn.Range("J2").Value = "fa02844e-356c-442b-89d9-021fb9e6330d.4bcf33a3301db87c2e28c06185ec766dc2fa2c70.zh-cn.xlf"

$wsZhCn.Range("A3").Value = "5ea61590-ac6e-4761-8d83-7c8a94ac5562.md"
Set-HyperlinkDisplay $wsZhCn '$A$3' "5ea61590-ac6e-4761-8d83-7c8a94ac5562.md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "5ea61590-ac6e-4761-8d83-7c8a94ac5562.0a86e44c14bffd7c6c4a405aeedc61d4493e0da9.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-04 06:52:23"
$wsZhCn.Range("I3").Value = "5ea61590-ac6e-4761-8d83-7c8a94ac5562.md"
Set-HyperlinkDisplay $wsZhCn '$I$3' "5ea61590-ac6e-4761-8d83-7c8a94ac5562.md"
$wsZhCn.Range("J3").Value = "5ea61590-ac6e-4761-8d83-7c8a94ac5562.0a86e44c14bffd7c6c4a405aeedc61d4493e0da9.zh-cn.xlf"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb47b9fc3bb9d46a7e3ed9c4965a617c6f0e9b97/e2e/5ea61590-ac6e-4761-8d83-7c8a94ac5562.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5cf4e2b670ca16205311b639d5287eb514fc0319/e2e/5ea61590-ac6e-4761-8d83-7c8a94ac5562.md."

# ---------------------------------------------------------------------------
# Sheet "de-de": mirrors the zh-cn changes (own dates / xlf file names).
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667

$wsDeDe.Range("A2").Value = "fa02844e-356c-442b-89d9-021fb9e6330d.md"
Set-HyperlinkDisplay $wsDeDe '$A$2' "fa02844e-356c-442b-89d9-021fb9e6330d.md"
$wsDeDe.Range("G2").Value = "fa02844e-356c-442b-89d9-021fb9e6330d.4bcf33a3301db87c2e28c06185ec766dc2fa2c70.de-de.xlf"
$wsDeDe.Range("I2").Value = "fa02844e-356c-442b-89d9-021fb9e6330d.md"
Set-HyperlinkDisplay $wsDeDe '$I$2' "fa02844e-356c-442b-89d9-021fb9e6330d.md"
$wsDeDe.Range("J2").Value = "fa02844e-356c-442b-89d9-021fb9e6330d.4bcf33a3301db87c2e28c06185ec766dc2fa2c70.de-de.xlf"

$wsDeDe.Range("A3").Value = "5ea61590-ac6e-4761-8d83-7c8a94ac5562.md"
Set-HyperlinkDisplay $wsDeDe '$A$3' "5ea61590-ac6e-4761-8d83-7c8a94ac5562.md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "5ea61590-ac6e-4761-8d83-7c8a94ac5562.0a86e44c14bffd7c6c4a405aeedc61d4493e0da9.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-04 06:52:28"
$wsDeDe.Range("I3").Value = "5ea61590-ac6e-4761-8d83-7c8a94ac5562.md"
Set-HyperlinkDisplay $wsDeDe '$I$3' "5ea61590-ac6e-4761-8d83-7c8a94ac5562.md"
$wsDeDe.Range("J3").Value = "5ea61590-ac6e-4761-8d83-7c8a94ac5562.0a86e44c14bffd7c6c4a405aeedc61d4493e0da9.de-de.xlf"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb47b9fc3bb9d46a7e3ed9c4965a617c6f0e9b97/e2e/5ea61590-ac6e-4761-8d83-7c8a94ac5562.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5cf4e2b670ca16205311b639d5287eb514fc0319/e2e/5ea61590-ac6e-4761-8d83-7c8a94ac5562.md."
